# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (1 -> 3) for rows 2-4 and
# propagate the recalculated expression / specificity / edge-weight values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.727457666666666
$ws.Range("H2").Value = 23.182373
$ws.Range("I2").Value = 0.1630271452636819
$ws.Range("J2").Value = 0.1630271452636819
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.276967
$ws.Range("N2").Value = 12.830901
$ws.Range("Q2").Value = 33.05008143423033
$ws.Range("R2").Value = 297.450732908073
$ws.Range("S2").Value = 0.1630271452636819
$ws.Range("T2").Value = 0.1630271452636819

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.42779933333333
$ws.Range("H3").Value = 112.283398
$ws.Range("I3").Value = 0.7896189849264272
$ws.Range("J3").Value = 0.7896189849264271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.276967
$ws.Range("N3").Value = 12.830901
$ws.Range("Q3").Value = 160.0774626312887
$ws.Range("R3").Value = 1440.697163681598
$ws.Range("S3").Value = 0.7896189849264272
$ws.Range("T3").Value = 0.7896189849264271

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.244565
$ws.Range("H4").Value = 6.733695
$ws.Range("I4").Value = 0.04735386980989085
$ws.Range("J4").Value = 0.04735386980989083
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.276967
$ws.Range("N4").Value = 12.830901
$ws.Range("Q4").Value = 9.599930434355
$ws.Range("R4").Value = 86.39937390919501
$ws.Range("S4").Value = 0.04735386980989085
$ws.Range("T4").Value = 0.04735386980989083
